$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.389.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +8.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.608.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9889"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3730"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3390"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.13"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.153"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9959"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.956"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.650"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9886"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001090"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.604.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06791"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +11.96%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.106"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.53%  "

$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.393.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.385"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.564"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +20.72%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.29%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.782.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.217"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.72%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.019"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +21.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9581"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08274"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.658"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.297"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.03%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06286"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.735"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.26%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.252"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02208"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2021"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.98%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6021"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.58%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9879"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.661"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5769"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.76%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.980"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06881"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.32%  "
